# Weekly refresh of the "Fruta / hortaliza" data: the price/volume rows
# are re-sorted by date, so the D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado) and P (Precio $/Kg)
# values move to new rows while the rest of each row's attributes
# (market, region, category, unit, origin, etc.) stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2..9, columns D, J, K, L, M, P
$data = @(
    @{ Row = 2;  D = 45245; J = 100; K = 9000;  L = 10000; M = 9500;  P = 528 },
    @{ Row = 3;  D = 45205; J = 200; K = 11000; L = 12000; M = 11500; P = 639 },
    @{ Row = 4;  D = 45175; J = 250; K = 11000; L = 12000; M = 11500; P = 639 },
    @{ Row = 5;  D = 45215; J = 200; K = 11000; L = 12000; M = 11500; P = 639 },
    @{ Row = 6;  D = 45092; J = 210; K = 10000; L = 11000; M = 10714; P = 595 },
    @{ Row = 7;  D = 44714; J = 80;  K = 9000;  L = 10000; M = 9500;  P = 528 },
    @{ Row = 8;  D = 44792; J = 160; K = 9000;  L = 10000; M = 9500;  P = 528 },
    @{ Row = 9;  D = 44804; J = 50;  K = 9500;  L = 10000; M = 9750;  P = 542 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("D$r").Value = $item.D
    $ws.Range("J$r").Value = $item.J
    $ws.Range("K$r").Value = $item.K
    $ws.Range("L$r").Value = $item.L
    $ws.Range("M$r").Value = $item.M
    $ws.Range("P$r").Value = $item.P
}
